# Applies the commit: inserts two new price records (rows) into the
# "Pimiento" sheet right before the existing row 230, shifting all
# subsequent rows down by two. The two new rows mirror the layout of
# the surrounding records (same market/region metadata) but carry
# their own date, volume and price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 230; all existing rows 230-324 move
# down to 232-326 (dimension grows from A1:R324 to A1:R326).
$ws.Rows.Item(230).Insert()
$ws.Rows.Item(230).Insert()

# ---- New row 230 : Zafiro rojo ----
$ws.Cells.Item(230, 1).Value = 5
$ws.Cells.Item(230, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(230, 3).Value = "Maule"
$ws.Cells.Item(230, 4).Value = 44466
$ws.Cells.Item(230, 5).Value = 7
$ws.Cells.Item(230, 6).Value = 100112002
$ws.Cells.Item(230, 7).Value = "Pimiento"
$ws.Cells.Item(230, 8).Value = "Zafiro rojo"
$ws.Cells.Item(230, 9).Value = "Primera"
$ws.Cells.Item(230, 10).Value = 200
$ws.Cells.Item(230, 11).Value = 40000
$ws.Cells.Item(230, 12).Value = 40000
$ws.Cells.Item(230, 13).Value = 40000
$ws.Cells.Item(230, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(230, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(230, 16).Value = 2667
$ws.Cells.Item(230, 17).Value = 15
$ws.Cells.Item(230, 18).Value = "Hortaliza"

# ---- New row 231 : Zafiro verde ----
$ws.Cells.Item(231, 1).Value = 5
$ws.Cells.Item(231, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(231, 3).Value = "Maule"
$ws.Cells.Item(231, 4).Value = 44466
$ws.Cells.Item(231, 5).Value = 7
$ws.Cells.Item(231, 6).Value = 100112002
$ws.Cells.Item(231, 7).Value = "Pimiento"
$ws.Cells.Item(231, 8).Value = "Zafiro verde"
$ws.Cells.Item(231, 9).Value = "Primera"
$ws.Cells.Item(231, 10).Value = 200
$ws.Cells.Item(231, 11).Value = 30000
$ws.Cells.Item(231, 12).Value = 30000
$ws.Cells.Item(231, 13).Value = 30000
$ws.Cells.Item(231, 14).Value = "`$/caja 15 kilos"
$ws.Cells.Item(231, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(231, 16).Value = 2000
$ws.Cells.Item(231, 17).Value = 15
$ws.Cells.Item(231, 18).Value = "Hortaliza"
